$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D ("Note"), shifting D:Q to F:S.
$ws.Range("D1:E1").EntireColumn.Insert()

# New header cells created by the insert (D1, E1 are currently blank).
$ws.Range("D1").Value = "Unnamed: 0.1.1"
$ws.Range("E1").Value = "Unnamed: 0.1.1.1"

# Fill the two new "Unnamed" index columns for every data row (same values as
# columns B/C - row index duplicated).
for ($r = 2; $r -le 8; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 4).Value = $idx   # column D
    $ws.Cells.Item($r, 5).Value = $idx   # column E
}

# Row 2 (ambient_temperature_system_failure.csv): MILOF now found an anomaly -
# update identified/overlap/params/time/best-param columns and fill the new
# best_paramKPar / best_paramBucket_index columns (R2, S2).
$ws.Range("M2").Value = "[3213, 3637, 6012]"
$ws.Range("N2").Value = 0.4
$ws.Range("O2").Value = "{'Numk': 17, 'KPar': 4, 'Bucket_index': 500}"
$ws.Range("P2").Value = 178.6288073339965
$ws.Range("Q2").Value = 17
$ws.Range("R2").Value = 4
$ws.Range("S2").Value = 500

# Row 3 (cpu_utilization_asg_misconfiguration.csv): same kind of update.
$ws.Range("M3").Value = "[16727, 17627, 17951]"
$ws.Range("N3").Value = 0.5
$ws.Range("O3").Value = "{'Numk': 23, 'KPar': 14, 'Bucket_index': 500}"
$ws.Range("P3").Value = 430.6199549960438
$ws.Range("Q3").Value = 23
$ws.Range("R3").Value = 14
$ws.Range("S3").Value = 500
